$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83, shifting existing rows 83.. down to 84..
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the new data record
$ws.Range("A83").Value = 4
$ws.Range("B83").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C83").Value = 'Los Lagos'
$ws.Range("D83").Value = 44518
$ws.Range("E83").Value = 10
$ws.Range("F83").Value = 100112039
$ws.Range("G83").Value = 'Ciboulette'
$ws.Range("H83").Value = 'Sin especificar'
$ws.Range("I83").Value = 'Primera'
$ws.Range("J83").Value = 80
$ws.Range("K83").Value = 2500
$ws.Range("L83").Value = 2500
$ws.Range("M83").Value = 2500
$ws.Range("N83").Value = '$/docena de atados'
$ws.Range("O83").Value = 'Región Metropolitana'
$ws.Range("P83").Value = 833
$ws.Range("Q83").Value = 3
$ws.Range("R83").Value = 'Hortaliza'

# D column uses a date-time number format (style index 2 in original) - match the format of surrounding date cells
$ws.Range("D83").NumberFormat = $ws.Range("D84").NumberFormat
